# "Export with no is_pref and no lev distance"
# Re-populate columns B (id) and C (speaker_variant) for rows 2-46 with the
# re-exported data (rows reordered/values updated per the new export), and
# clear column D (is_prefered) for every data row since the new export no
# longer marks any row as preferred.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = "#griet"
$ws.Cells.Item(2, 3).Value = "griet"
$ws.Cells.Item(3, 2).Value = "#alph"
$ws.Cells.Item(3, 3).Value = "Alph"
$ws.Cells.Item(4, 2).Value = "#elys"
$ws.Cells.Item(4, 3).Value = "Elys"
$ws.Cells.Item(5, 2).Value = "#rod"
$ws.Cells.Item(5, 3).Value = "Rod"
$ws.Cells.Item(6, 2).Value = "#koning"
$ws.Cells.Item(6, 3).Value = "Koning"
$ws.Cells.Item(7, 2).Value = "#almijn"
$ws.Cells.Item(7, 3).Value = "Almijn"
$ws.Cells.Item(8, 2).Value = "#elisabet"
$ws.Cells.Item(8, 3).Value = "Elisabet"
$ws.Cells.Item(9, 2).Value = "#nieu.-h"
$ws.Cells.Item(9, 3).Value = "Nieu. H"
$ws.Cells.Item(10, 2).Value = "#roderick"
$ws.Cells.Item(10, 3).Value = "Roderick"
$ws.Cells.Item(11, 2).Value = "#oronda"
$ws.Cells.Item(11, 3).Value = "Oronda"
$ws.Cells.Item(12, 2).Value = "#nieu.-hae"
$ws.Cells.Item(12, 3).Value = "Nieu. hae"
$ws.Cells.Item(13, 2).Value = "#griet"
$ws.Cells.Item(13, 3).Value = "Griet"
$ws.Cells.Item(14, 2).Value = "#coningh"
$ws.Cells.Item(14, 3).Value = "Coningh"
$ws.Cells.Item(15, 2).Value = "#haalnaa"
$ws.Cells.Item(15, 3).Value = "Haalnaa"
$ws.Cells.Item(16, 2).Value = "#alphonsus"
$ws.Cells.Item(16, 3).Value = "Alphonsus"
$ws.Cells.Item(17, 2).Value = "#elisab"
$ws.Cells.Item(17, 3).Value = "Elisab"
$ws.Cells.Item(18, 2).Value = "#alphon"
$ws.Cells.Item(18, 3).Value = "Alphon"
$ws.Cells.Item(19, 2).Value = "#alph.-b"
$ws.Cells.Item(19, 3).Value = "Alph. b"
$ws.Cells.Item(20, 2).Value = "#griet-s"
$ws.Cells.Item(20, 3).Value = "Griet S"
$ws.Cells.Item(21, 2).Value = "#gerald"
$ws.Cells.Item(21, 3).Value = "Gerald"
$ws.Cells.Item(22, 2).Value = "#geerald"
$ws.Cells.Item(22, 3).Value = "Geerald"
$ws.Cells.Item(23, 2).Value = "#elisabeth"
$ws.Cells.Item(23, 3).Value = "Elisabeth"
$ws.Cells.Item(24, 2).Value = "#elysabet"
$ws.Cells.Item(24, 3).Value = "Elysabet"
$ws.Cells.Item(25, 2).Value = "#elijsab"
$ws.Cells.Item(25, 3).Value = "Elijsab"
$ws.Cells.Item(26, 2).Value = "#koningh"
$ws.Cells.Item(26, 3).Value = "Koningh"
$ws.Cells.Item(27, 2).Value = "#rodrick"
$ws.Cells.Item(27, 3).Value = "Rodrick"
$ws.Cells.Item(28, 2).Value = "#bode"
$ws.Cells.Item(28, 3).Value = "Bode"
$ws.Cells.Item(29, 2).Value = "#elysab"
$ws.Cells.Item(29, 3).Value = "Elysab"
$ws.Cells.Item(30, 2).Value = "#alphons"
$ws.Cells.Item(30, 3).Value = "Alphons"
$ws.Cells.Item(31, 2).Value = "#roderic"
$ws.Cells.Item(31, 3).Value = "Roderic"
$ws.Cells.Item(32, 2).Value = "#elysa"
$ws.Cells.Item(32, 3).Value = "Elysa"
$ws.Cells.Item(33, 2).Value = "#rodrick,"
$ws.Cells.Item(33, 3).Value = "Rodrick,"
$ws.Cells.Item(34, 2).Value = "#orond"
$ws.Cells.Item(34, 3).Value = "Orond"
$ws.Cells.Item(35, 2).Value = "#pagie"
$ws.Cells.Item(35, 3).Value = "Pagie"
$ws.Cells.Item(36, 2).Value = "#kamenier"
$ws.Cells.Item(36, 3).Value = "Kamenier"
$ws.Cells.Item(37, 2).Value = "#gerald"
$ws.Cells.Item(37, 3).Value = "gerald"
$ws.Cells.Item(38, 2).Value = "#nieu.-ha"
$ws.Cells.Item(38, 3).Value = "Nieu. ha"
$ws.Cells.Item(39, 2).Value = "#elijsabet"
$ws.Cells.Item(39, 3).Value = "Elijsabet"
$ws.Cells.Item(40, 2).Value = "#rodd"
$ws.Cells.Item(40, 3).Value = "Rodd"
$ws.Cells.Item(41, 2).Value = "#gerald,"
$ws.Cells.Item(41, 3).Value = "Gerald,"
$ws.Cells.Item(42, 2).Value = "#rodde"
$ws.Cells.Item(42, 3).Value = "Rodde"
$ws.Cells.Item(43, 2).Value = "#haalna"
$ws.Cells.Item(43, 3).Value = "Haalna"
$ws.Cells.Item(44, 2).Value = "#nieuw.-h"
$ws.Cells.Item(44, 3).Value = "Nieuw. H"
$ws.Cells.Item(45, 2).Value = "#alphonse"
$ws.Cells.Item(45, 3).Value = "Alphonse"
$ws.Cells.Item(46, 2).Value = "#rodderick"
$ws.Cells.Item(46, 3).Value = "Rodderick"

# Clear D2:D20 (remove is_prefered "x" marks)
$ws.Range("D2:D20").ClearContents()
